$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'245.88"
$ws.Range("D3").Formula = "'22.10"
$ws.Range("D5").Formula = "'0.05873"
$ws.Range("D6").Formula = "'3.382"
$ws.Range("D7").Formula = "'6.373"
$ws.Range("D8").Formula = "'0.8172"
$ws.Range("D9").Formula = "'0.9586"
$ws.Range("D10").Formula = "'0.1420"
$ws.Range("D11").Formula = "'0.03593"
$ws.Range("D12").Formula = "'0.07333"
$ws.Range("D13").Formula = "'0.03040"
$ws.Range("D14").Formula = "'4.425"
$ws.Range("D15").Formula = "'0.09390"
$ws.Range("D16").Formula = "'0.001588"
$ws.Range("D17").Formula = "'0.04807"
$ws.Range("D18").Formula = "'0.0005902"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Formula = "'0.006105"
$ws.Range("D20").Formula = "'0.004077"
$ws.Range("D21").Formula = "'0.0009877"
$ws.Range("D22").Formula = "'0.00009711"
$ws.Range("D23").Formula = "'3.688"
$ws.Range("D25").Formula = "'0.3263"
$ws.Range("D27").Formula = "'0.0002472"
$ws.Range("D40").Formula = "'0.03858"
$ws.Range("D41").Formula = "'0.006587"
$ws.Range("D42").Formula = "'0.1074"
$ws.Range("D43").Formula = "'0.002443"
$ws.Range("D44").Formula = "'0.005897"
$ws.Range("D45").Formula = "'0.00005668"
$ws.Range("D47").Formula = "'0.7752"
$ws.Range("D48").Formula = "'0.05425"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
